$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("20210920-20%")

# The row-number helper column (A11:A310) used a shared formula "ROW()-9"
# which produced an off-by-one count (row 11 -> 2 instead of 1). Fix the
# formula (and thus the cached values) so that row 11 -> 1, row 12 -> 2, etc.
$ws.Range("A11:A310").Formula = "=ROW()-10"
